$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 ("На всякий случай"): update values in B7 and C7
$ws.Range("B7").Value = -34210
$ws.Range("C7").Value = 34912

# Row 8 ("Еда"): update values in B8 and C8
$ws.Range("B8").Value = 2520
$ws.Range("C8").Value = 980
